{"js": "// Replace the date line and the 26 multiplication problems with their\n// new values. Each original string is unique in the document, so a\n// simple search-and-replace per pair is sufficient and order-independent.\nconst replacements = [\n  [\"2025-10-14 Tuesday\", \"2025-10-15 Wednesday\"],\n  [\"25\u00d725=\", \"44\u00d780=\"],\n  [\"15\u00d799=\", \"62\u00d745=\"],\n  [\"57\u00d770=\", \"95\u00d788=\"],\n  [\"52\u00d776=\", \"83\u00d787=\"],\n  [\"47\u00d777=\", \"31\u00d754=\"],\n  [\"91\u00d747=\", \"60\u00d769=\"],\n  [\"78\u00d780=\", \"64\u00d728=\"],\n  [\"72\u00d763=\", \"11\u00d787=\"],\n  [\"35\u00d766=\", \"42\u00d791=\"],\n  [\"66\u00d746=\", \"11\u00d737=\"],\n  [\"96\u00d794=\", \"50\u00d718=\"],\n  [\"79\u00d729=\", \"62\u00d767=\"],\n  [\"54\u00d736=\", \"79\u00d718=\"],\n  [\"97\u00d780=\", \"74\u00d726=\"],\n  [\"99\u00d743=\", \"61\u00d775=\"],\n  [\"11\u00d749=\", \"77\u00d741=\"],\n  [\"94\u00d717=\", \"86\u00d759=\"],\n  [\"27\u00d772=\", \"77\u00d762=\"],\n  [\"72\u00d761=\", \"42\u00d775=\"],\n  [\"84\u00d788=\", \"93\u00d721=\"],\n  [\"33\u00d780=\", \"42\u00d777=\"],\n  [\"87\u00d731=\", \"74\u00d731=\"],\n  [\"72\u00d784=\", \"17\u00d766=\"],\n  [\"77\u00d756=\", \"15\u00d797=\"],\n  [\"37\u00d738=\", \"93\u00d730=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 26 multiplication problems with their\n# new values. Each original string is unique in the document, so a\n# Find/Replace pass per pair (restricted to one replacement via\n# wdReplaceOne) is sufficient and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-14 Tuesday\", \"2025-10-15 Wednesday\"),\n    @(\"25\u00d725=\", \"44\u00d780=\"),\n    @(\"15\u00d799=\", \"62\u00d745=\"),\n    @(\"57\u00d770=\", \"95\u00d788=\"),\n    @(\"52\u00d776=\", \"83\u00d787=\"),\n    @(\"47\u00d777=\", \"31\u00d754=\"),\n    @(\"91\u00d747=\", \"60\u00d769=\"),\n    @(\"78\u00d780=\", \"64\u00d728=\"),\n    @(\"72\u00d763=\", \"11\u00d787=\"),\n    @(\"35\u00d766=\", \"42\u00d791=\"),\n    @(\"66\u00d746=\", \"11\u00d737=\"),\n    @(\"96\u00d794=\", \"50\u00d718=\"),\n    @(\"79\u00d729=\", \"62\u00d767=\"),\n    @(\"54\u00d736=\", \"79\u00d718=\"),\n    @(\"97\u00d780=\", \"74\u00d726=\"),\n    @(\"99\u00d743=\", \"61\u00d775=\"),\n    @(\"11\u00d749=\", \"77\u00d741=\"),\n    @(\"94\u00d717=\", \"86\u00d759=\"),\n    @(\"27\u00d772=\", \"77\u00d762=\"),\n    @(\"72\u00d761=\", \"42\u00d775=\"),\n    @(\"84\u00d788=\", \"93\u00d721=\"),\n    @(\"33\u00d780=\", \"42\u00d777=\"),\n    @(\"87\u00d731=\", \"74\u00d731=\"),\n    @(\"72\u00d784=\", \"17\u00d766=\"),\n    @(\"77\u00d756=\", \"15\u00d797=\"),\n    @(\"37\u00d738=\", \"93\u00d730=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
